# Working on reviewer comments: add a "Replication 4 sets" results table
# to the "titan" worksheet (rows 64-71) of the scaling workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Write the text labels first, in the exact order the original author
# typed them, so the shared-string table is built up in the same sequence
# (35=Replication 4 sets, 36=Run, 37=Set-AllReduce, 38=Set-AllToAll,
#  39=Block-AllReduce, 40=Block-AllToAll, 41=MC Min, 42=MC Ave, 43=MC Max,
#  44=Sum Min, 45=Sum Max, 46=Sum Ave, 47=Min, 48=Ave, 49=Max).
$ws.Cells.Item(64, 1).Value = "Replication 4 sets"
$ws.Cells.Item(65, 1).Value = "Run"
$ws.Cells.Item(66, 1).Value = "Set-AllReduce"
$ws.Cells.Item(68, 1).Value = "Set-AllToAll"
$ws.Cells.Item(67, 1).Value = "Block-AllReduce"
$ws.Cells.Item(69, 1).Value = "Block-AllToAll"
$ws.Cells.Item(65, 2).Value = "MC Min"
$ws.Cells.Item(65, 3).Value = "MC Ave"
$ws.Cells.Item(65, 4).Value = "MC Max"
$ws.Cells.Item(65, 5).Value = "Sum Min"
$ws.Cells.Item(65, 7).Value = "Sum Max"
$ws.Cells.Item(65, 6).Value = "Sum Ave"
$ws.Cells.Item(65, 8).Value = "Min"
$ws.Cells.Item(65, 9).Value = "Ave"
$ws.Cells.Item(65, 10).Value = "Max"

# --- Numeric data for the two raw-measurement rows -------------------------
$ws.Cells.Item(66, 2).Value = 3.5470000000000002
$ws.Cells.Item(66, 3).Value = 3.5653999999999999
$ws.Cells.Item(66, 4).Value = 3.5842999999999998
$ws.Cells.Item(66, 5).Value = 0.041291
$ws.Cells.Item(66, 6).Value = 0.13708000000000001
$ws.Cells.Item(66, 7).Value = 0.19917000000000001

$ws.Cells.Item(67, 2).Value = 3.7547999999999999
$ws.Cells.Item(67, 3).Value = 3.7799
$ws.Cells.Item(67, 4).Value = 3.8121
$ws.Cells.Item(67, 5).Value = 0.031575
$ws.Cells.Item(67, 6).Value = 0.088169
$ws.Cells.Item(67, 7).Value = 0.15418999999999999

# --- Formulas --------------------------------------------------------------
# Row 66 gets plain (non-shared) formulas.
$ws.Range("H66").Formula = "=B66+E66"
$ws.Range("I66").Formula = "=C66+F66"
$ws.Range("J66").Formula = "=D66+G66"

# Rows 67-69 share one formula per column (H/I/J), matching how Excel
# records a fill-down over a multi-row range.
$ws.Range("H67:H69").Formula = "=B67+E67"
$ws.Range("I67:I69").Formula = "=C67+F67"
$ws.Range("J67:J69").Formula = "=D67+G67"

# --- Formatting --------------------------------------------------------------
# Section title: underlined, no fill/border (matches the workbook's other
# plain underlined labels).
$ws.Cells.Item(64, 1).Font.Underline = $true

# Header row (A65:J65): bold + right aligned.
$ws.Range("A65:J65").Font.Bold = $true
$ws.Range("A65:J65").HorizontalAlignment = -4152   # xlRight

# Row-label cells A66/A67 ("Set-AllReduce"/"Block-AllReduce"): bold + right aligned.
$ws.Range("A66").Font.Bold = $true
$ws.Range("A66").HorizontalAlignment = -4152
$ws.Range("A67").Font.Bold = $true
$ws.Range("A67").HorizontalAlignment = -4152

# Numeric rows 66/67 (B:J): right aligned, regular weight.
$ws.Range("B66:J66").HorizontalAlignment = -4152
$ws.Range("B67:J67").HorizontalAlignment = -4152

# A68/A69/A71 labels: bold (same weight as other sub-totals in the sheet).
$ws.Range("A68").Font.Bold = $true
$ws.Range("A69").Font.Bold = $true
$ws.Range("A71").Font.Bold = $true

# --- Column A width (12.5 -> 15 characters) --------------------------------
$ws.Columns.Item(1).ColumnWidth = 14.17

# --- Sheet view: scroll position + active selection -------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("A65:J67").Select()
